# Update EnemyData worksheet to match latest table format.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EnemyData")

# Update the two marker strings in row 5 (shared strings table entries).
$ws.Range("A5").Value = "skip"
$ws.Range("B5:K5").Value = "both"

# Move the active selection from E14 to A6.
$ws.Range("A6").Select()
